# "Changes of 23rd June Final"
# Updates the FedEx rate-verification sheet: column E (ActualRate) gets the
# newly observed actual rates, and column F (Result) flips from PASS to
# FAIL for every data row (2-31) now that actual != expected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$actualRates = @(
    "$18.40", "$19.73", "$27.61", "$38.85", "$41.36", "$63.39", "$74.80", "$210.32",
    "$6.00", "$6.00", "$24.00", "$36.00", "$6.00", "$63.00", "$6.18", "$13.55",
    "$18.91", "$29.42", "$35.13", "$48.33", "$37.82", "$48.33", "$56.73", "$79.85",
    "$105.06", "$105.06", "$148.16", "$287.20", "$476.40", "$273.68"
)

$startRow = 2
for ($i = 0; $i -lt $actualRates.Length; $i++) {
    $row = $startRow + $i

    # Column E ("ActualRate") must stay a plain text string like "$18.40",
    # not get auto-coerced into a currency number by Excel's parser. Format
    # the cell as Text first, write the value, then restore the cell style
    # so no numeric/currency formatting sticks to it.
    $eCell = $ws.Range("E$row")
    $eCell.NumberFormat = "@"
    $eCell.Value = $actualRates[$i]
    $eCell.Style = "Normal"

    # Column F ("Result") always flips to FAIL now.
    $ws.Range("F$row").Value = "FAIL"
}
